$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuoteOptionPageData")
$insured = $wb.Worksheets.Item("InsuredPageData")

# ---------------------------------------------------------------------------
# New block on QuoteOptionPageData: "testConfirmAndLockQuoteOption" (rows 11-13)
# Built by copying the styling of the existing "testDeleteQuoteOption" block
# (rows 6-8) which has the same row layout, minus the "website" column.
# ---------------------------------------------------------------------------

# Row 11: section title (copies border/fill styling used by the other titles)
$ws.Range("A6").Copy($ws.Range("A11"))
$ws.Range("A11").Value = "testConfirmAndLockQuoteOption"

# Row 12: header labels (copy row 7's header, skipping the "website" column D)
$ws.Range("A7:C7").Copy($ws.Range("A12"))
$ws.Range("E7:O7").Copy($ws.Range("D12"))

# Row 13: data values (copy row 8's data, skipping the "website" column D)
$ws.Range("A8:C8").Copy($ws.Range("A13"))
$ws.Range("E8:O8").Copy($ws.Range("D13"))

# Now overwrite the values that differ from the copied template
$ws.Range("A13").Value = "Y"
$ws.Range("B13").Value = "NetGuard® Plus"
$ws.Range("C13").Value = "Active"
$ws.Range("D13").Value = 20217
$ws.Range("E13").Value = 237
# F13 keeps a plain (non-alignment) border style, unlike the other numeric cells
$ws.Range("C8").Copy($ws.Range("F13"))
$ws.Range("F13").Value = 8006
$ws.Range("G13").Value = 173
$ws.Range("H13").Value = 237
$ws.Range("I13").Value = "Business to Business"
$ws.Range("J13").Value = 1000000
$ws.Range("K13").Value = 445
$ws.Range("L13").Value = "$ 2MM"
$ws.Range("M13").Value = "$ 5MM"
$ws.Range("N13").Value = 12000

# Column sizing for the widened/new columns used by the new block
$ws.Columns.Item(1).ColumnWidth = 28.6640625
$ws.Columns.Item(4).ColumnWidth = 27
$ws.Columns.Item(6).ColumnWidth = 13.33203125
$ws.Columns.Item(7).ColumnWidth = 13.33203125
$ws.Columns.Item(8).ColumnWidth = 22.6640625
$ws.Columns.Item(9).ColumnWidth = 24.109375

# ---------------------------------------------------------------------------
# Move the active tab / selection from InsuredPageData to QuoteOptionPageData
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("B15").Select()
